# Weekly update: insert a new price-record row above the current row 53
# (shifting all the existing Arveja Verde records for
# "Comercializadora del Agro de Limarí" down by one row) and populate the
# new row with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 53; everything currently at 53:79
# shifts down to 54:80 (Excel also copies formatting, e.g. the date
# number-format on column D, from the row above).
$ws.Rows("53:53").Insert()

# Populate the newly inserted row 53 with the new weekly record.
$ws.Cells.Item(53, 1).Value = 2
$ws.Cells.Item(53, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(53, 3).Value = "Coquimbo"
$ws.Cells.Item(53, 4).Value = 44846
$ws.Cells.Item(53, 5).Value = 4
$ws.Cells.Item(53, 6).Value = 100112022
$ws.Cells.Item(53, 7).Value = "Arveja Verde"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 488
$ws.Cells.Item(53, 11).Value = 23000
$ws.Cells.Item(53, 12).Value = 24000
$ws.Cells.Item(53, 13).Value = 23426
$ws.Cells.Item(53, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(53, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(53, 16).Value = 937
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 18).Value = "Hortaliza"
